# Daily attendance processing - 2025-10-29 16:27:10
# For every data row in the "Recorded By" column (G), swap the first two
# comma-separated entries while leaving any additional entries (e.g. a
# trailing "system" marker) in place. Cells with a single entry (no comma)
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -eq $null -or $text -eq "") {
        continue
    }

    $parts = $text.Split(",")

    if ($parts.Length -ge 2) {
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }

        $tmp = $parts[0]
        $parts[0] = $parts[1]
        $parts[1] = $tmp

        $cell.Value = [string]::Join(", ", $parts)
    }
}
